$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 3
$ws.Range("I3").Value = 2.38
$ws.Range("R3").Value = 1.7
$ws.Range("S3").Value = 2.05
$ws.Range("G4").Value = 8
$ws.Range("K4").Value = 21
$ws.Range("L4").Value = 1.14
$ws.Range("M4").Value = 5.5
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 2.05
$ws.Range("AG4").Value = 8.5
$ws.Range("G7").Value = 1.65
$ws.Range("J7").Value = 1.04
$ws.Range("K7").Value = 13
$ws.Range("L7").Value = 1.22
$ws.Range("M7").Value = 4
$ws.Range("AC7").Value = 51
$ws.Range("G8").Value = 1.36
$ws.Range("H8").Value = 4.75
$ws.Range("K8").Value = 19
$ws.Range("L8").Value = 1.14
$ws.Range("M8").Value = 5.5
$ws.Range("N8").Value = 1.48
$ws.Range("O8").Value = 2.6
$ws.Range("V8").Value = 9
$ws.Range("X8").Value = 11
$ws.Range("Y8").Value = 21
$ws.Range("Z8").Value = 19
$ws.Range("AA8").Value = 9.5
$ws.Range("G9").Value = 1.91
$ws.Range("H9").Value = 3.2
$ws.Range("AA9").Value = 6
$ws.Range("AF9").Value = 21
$ws.Range("AH9").Value = 51
$ws.Range("G10").Value = 2.2
$ws.Range("I10").Value = 2.82
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3.65
$ws.Range("L11").Value = 1.39
$ws.Range("M11").Value = 2.55
$ws.Range("N11").Value = 2.12
$ws.Range("O11").Value = 1.57
$ws.Range("P11").Value = 1.47
$ws.Range("Q11").Value = 2.32
$ws.Range("R11").Value = 1.93
$ws.Range("S11").Value = 1.7
$ws.Range("T11").Value = 6
$ws.Range("U11").Value = 8.5
$ws.Range("V11").Value = 9
$ws.Range("X11").Value = 18
$ws.Range("Y11").Value = 35
$ws.Range("Z11").Value = 7.8
$ws.Range("AB11").Value = 17.5
$ws.Range("AC11").Value = 100
$ws.Range("AE11").Value = 9
$ws.Range("AF11").Value = 18
$ws.Range("AG11").Value = 13
$ws.Range("AI11").Value = 37
$ws.Range("AJ11").Value = 50
$ws.Range("G13").Value = 2.1
$ws.Range("G17").Value = 2.2
$ws.Range("I17").Value = 3.7
$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 1.75
$ws.Range("U17").Value = 9.5
$ws.Range("Z17").Value = 7
$ws.Range("AB17").Value = 17
$ws.Range("AD17").Value = 401
$ws.Range("AI17").Value = 34
$ws.Range("G19").Value = 2
$ws.Range("T19").Value = 8.5
$ws.Range("U19").Value = 10
$ws.Range("W19").Value = 17
$ws.Range("Z19").Value = 13
$ws.Range("G23").Value = 3.25
$ws.Range("L23").Value = 1.33
$ws.Range("M23").Value = 3.25
$ws.Range("N23").Value = 2.05
$ws.Range("O23").Value = 1.75
$ws.Range("R23").Value = 1.83
$ws.Range("S23").Value = 1.83
$ws.Range("T23").Value = 9.5
$ws.Range("Z23").Value = 9.5
$ws.Range("AI23").Value = 19
$ws.Range("H30").Value = 3.9
$ws.Range("I30").Value = 4.1
$ws.Range("H31").Value = 4.1
$ws.Range("I31").Value = 6.25
$ws.Range("K31").Value = 12
$ws.Range("L31").Value = 1.22
$ws.Range("M31").Value = 4
$ws.Range("T31").Value = 7.5
$ws.Range("AD31").Value = 251
$ws.Range("AE31").Value = 17
$ws.Range("AF31").Value = 34
$ws.Range("N32").Value = 1.73
$ws.Range("O32").Value = 2.08
$ws.Range("H33").Value = 4.1
$ws.Range("I33").Value = 1.53
$ws.Range("N33").Value = 1.84
$ws.Range("O33").Value = 1.89
$ws.Range("AJ33").Value = 29
